# feat: add 2022-Q4 data
# ------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert the new 2022-Q4 row at the top of
#    the data table, shifting every other row down by one and
#    re-creating the row that falls off the bottom (2020-Q4).
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# Clone row 8's look into a brand-new row 9 (keeps the bordered/bold
# style on column A identical to the rest of the table), then rewrite
# every row's values so row 2 becomes 2022-Q4 and the rest shift down.
$total.Range("A8:D8").Copy($total.Range("A9:D9"))

$total.Range("A9").Value = 7
$total.Range("B9").Value = "2020-Q4"
$total.Range("C9").Value = 2
$total.Range("D9").Value = 0.02

$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 2
$total.Range("D8").Value = 0.07000000000000001

$total.Range("B7").Value = "2021-Q3"
$total.Range("C7").Value = 5
$total.Range("D7").Value = 6.35

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 5.97

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 9.09

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 0.57

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.03

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 5.48

# ------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e.
#    right before "2022-Q3"). Duplicating the existing "2022-Q3" sheet
#    keeps every sheet-level setting (outline props, page margins,
#    header/index-column styles) identical instead of inventing new
#    style or sheetPr records for a freshly blank sheet.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The duplicated sheet only has 2 fund rows (rows 2-3); the new quarter
# needs 5, so stamp out 3 more rows using row 3's formatting as the
# template (keeps the same styles/number formats for every column).
$q4.Range("A3:H3").Copy($q4.Range("A4:H4"))
$q4.Range("A3:H3").Copy($q4.Range("A5:H5"))
$q4.Range("A3:H3").Copy($q4.Range("A6:H6"))

# -- Header row --
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# -- Index column (A): plain integers 0..4 --
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1
$q4.Range("A4").Value = 2
$q4.Range("A5").Value = 3
$q4.Range("A6").Value = 4

# Columns B-G store numeric-looking values as literal text in the
# source data, so force Text format before assigning or Excel would
# silently coerce them (and strip leading zeros from fund codes).
$q4.Range("B2:G6").NumberFormat = "@"

$q4.Range("B2").Value = "005669"
$q4.Range("C2").Value = "前海开源公用事业行业股票"
$q4.Range("D2").Value = "168.84"
$q4.Range("E2").Value = "94.23"
$q4.Range("F2").Value = "3.20"
$q4.Range("G2").Value = "5.4029"
$q4.Range("H2").Value = 10

$q4.Range("B3").Value = "006923"
$q4.Range("C3").Value = "前海开源沪港深非周期性行业股票A"
$q4.Range("D3").Value = "0.28"
$q4.Range("E3").Value = "90.65"
$q4.Range("F3").Value = "7.24"
$q4.Range("G3").Value = "0.0203"
$q4.Range("H3").Value = 4

$q4.Range("B4").Value = "010540"
$q4.Range("C4").Value = "浙商智多金稳健一年持有期混合C"
$q4.Range("D4").Value = "1.37"
$q4.Range("E4").Value = "25.01"
$q4.Range("F4").Value = "1.47"
$q4.Range("G4").Value = "0.0201"
$q4.Range("H4").Value = 6

$q4.Range("B5").Value = "010539"
$q4.Range("C5").Value = "浙商智多金稳健一年持有期混合A"
$q4.Range("D5").Value = "1.27"
$q4.Range("E5").Value = "25.01"
$q4.Range("F5").Value = "1.47"
$q4.Range("G5").Value = "0.0187"
$q4.Range("H5").Value = 6

$q4.Range("B6").Value = "006924"
$q4.Range("C6").Value = "前海开源沪港深非周期性行业股票C"
$q4.Range("D6").Value = "0.24"
$q4.Range("E6").Value = "90.65"
$q4.Range("F6").Value = "7.24"
$q4.Range("G6").Value = "0.0174"
$q4.Range("H6").Value = 4

$total.Select()
$total.Range("A1").Select()
